$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "46.181.86"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -1.29%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.364.50"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +2.56%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$ws.Range("E4").Value = "  -0.22%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "301.73"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.92%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "100.33"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -1.45%  "

$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("E8").Value = "  -0.04%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.514"
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "34.57"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -6.05%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0800"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("E12").Value = "  -3.04%  "

$ws.Range("E13").Value = "  -0.36%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.723.41"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +2.49%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.359.81"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +2.22%  "

$ws.Range("E16").Value = "  -0.68%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "13.66"
$cell.ClearFormats()

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "46.105.66"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -1.38%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.82"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -2.60%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0968"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +2.61%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.05"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -1.73%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "67.75"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +1.03%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "246.04"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -0.69%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.82"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -3.86%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -0.10%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.92"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -2.71%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "40.03"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -6.21%  "

$ws.Range("E28").Value = "  -2.69%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.82"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -1.21%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.80"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +22.37%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "21.04"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +4.02%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.81"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +7.46%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.53"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -3.77%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "146.41"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -0.33%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.0778"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -2.69%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.113"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -0.30%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.91"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +6.59%  "

$ws.Range("E38").Value = "  -2.29%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "15.08"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -4.60%  "

$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("E42").Value = "  -6.41%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.904.40"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +3.38%  "

$ws.Range("E44").Value = "  -0.05%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "92.74"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +3.42%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.82"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -8.28%  "

$ws.Range("E47").Value = "  -5.40%  "

$ws.Range("E48").Value = "  +4.53%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "98.15"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +0.68%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.596.56"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +2.41%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "69.48"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -7.74%  "
